# ExcelTest.xlsx - add AllProductsPage tests
# 1) Update the absPath hint, rename sheet1, add two new sheets
# 2) Populate the two new sheets with test data
# 3) Fix up selections/active sheet to match the authored state

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# ---------------------------------------------------------------------------
# Sheets: rename existing "validateLogin" -> "verifyLogin", add two new sheets
# named "verifyAddAndRemoveButtons" and "verifyItemSort" after it.
# An extra throwaway sheet is added+removed first purely so the sheetId
# counter advances the same way it did for the original author (sheetId 3/4
# instead of 2/3).
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item(1)
$wsLogin.Name = "verifyLogin"

$wsThrowaway = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsThrowaway.Name = "zzThrowaway"

$wsButtons = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsButtons.Name = "verifyAddAndRemoveButtons"

$wsSort = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsSort.Name = "verifyItemSort"

$wb.Worksheets.Item("zzThrowaway").Delete()

# Re-fetch live references to the sheets (indices shift after Delete)
$wsLogin = $wb.Worksheets.Item("verifyLogin")
$wsButtons = $wb.Worksheets.Item("verifyAddAndRemoveButtons")
$wsSort = $wb.Worksheets.Item("verifyItemSort")

# ---------------------------------------------------------------------------
# verifyAddAndRemoveButtons sheet data
# ---------------------------------------------------------------------------
$wsButtons.Range("A1").Value = "itemsToAdd"
$wsButtons.Range("A2").Value = "Sauce Labs Backpack"
$wsButtons.Range("A3").Value = "Sauce Labs Bolt T-Shirt, Sauce Labs Backpack, Test.allTheThings() T-Shirt (Red)"
$wsButtons.Range("A4").Value = "Sauce Labs Backpack, Sauce Labs Bike Light, Sauce Labs Bolt T-Shirt, Sauce Labs Fleece Jacket, Sauce Labs Onesie, Test.allTheThings() T-Shirt (Red)"

$wsButtons.Range("D3").Value = "'1"
$wsButtons.Range("C4").Value = "'3"
$wsButtons.Range("C2").Value = "'6"

$wsButtons.Range("B1").Value = "itemsToRemove"
$wsButtons.Range("D2").Value = "'0"
$wsButtons.Range("B3").Value = "Sauce Labs Bolt T-Shirt, Test.allTheThings() T-Shirt (Red)"

$wsButtons.Range("C1").Value = "addButtonCount"
$wsButtons.Range("D1").Value = "removeButtonCount"
$wsButtons.Range("C3").Value = "'5"
$wsButtons.Range("B4").Value = "Sauce Labs Bike Light, Sauce Labs Bolt T-Shirt, Test.allTheThings() T-Shirt (Red)"

# ---------------------------------------------------------------------------
# verifyItemSort sheet data
# ---------------------------------------------------------------------------
$wsSort.Range("A1").Value = "sortType"
$wsSort.Range("B1").Value = "firstItem"
$wsSort.Range("C1").Value = "lastItem"
$wsSort.Range("A2").Value = "Name (Z to A)"
$wsSort.Range("A3").Value = "Price (low to high)"
$wsSort.Range("A4").Value = "Price (high to low)"
$wsSort.Range("A5").Value = "Name (A to Z)"
$wsSort.Range("B2").Value = "Test.allTheThings() T-Shirt (Red)"
$wsSort.Range("B3").Value = "Sauce Labs Onesie"
$wsSort.Range("C3").Value = "Sauce Labs Fleece Jacket"

# cartCount column added last (after the sort sheet was filled in)
$wsButtons.Range("E1").Value = "cartCount"

# ---------------------------------------------------------------------------
# Remaining cells that just reuse already-introduced values
# ---------------------------------------------------------------------------
$wsButtons.Range("B2").Value = "Sauce Labs Backpack"
$wsButtons.Range("E2").Value = "'0"
$wsButtons.Range("E3").Value = "'1"
$wsButtons.Range("D4").Value = "'3"
$wsButtons.Range("E4").Value = "'3"

$wsSort.Range("C2").Value = "Sauce Labs Backpack"
$wsSort.Range("B4").Value = "Sauce Labs Fleece Jacket"
$wsSort.Range("C4").Value = "Sauce Labs Onesie"
$wsSort.Range("B5").Value = "Sauce Labs Backpack"
$wsSort.Range("C5").Value = "Test.allTheThings() T-Shirt (Red)"

# verifyItemSort was printed in portrait orientation
$wsSort.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Column widths (bestFit look) on the new sheets - set to (approximately)
# the widths Excel's own AutoFit produced for this content/font.
# ---------------------------------------------------------------------------
$wsButtons.Columns.Item(1).ColumnWidth = 131.45182291666666
$wsButtons.Columns.Item(2).ColumnWidth = 71.30729166666667
$wsButtons.Columns.Item(3).ColumnWidth = 14.877604166666666
$wsButtons.Columns.Item(4).ColumnWidth = 18.592447916666664

$wsSort.Columns.Item(1).ColumnWidth = 16.592447916666664
$wsSort.Columns.Item(2).ColumnWidth = 18.307291666666664
$wsSort.Columns.Item(3).ColumnWidth = 29.307291666666664

# ---------------------------------------------------------------------------
# Selections / active sheet to match the saved authoring state
# ---------------------------------------------------------------------------
$wsLogin.Activate()
$wsLogin.Range("C8").Select()

$wsSort.Activate()
$wsSort.Range("C1").Select()

$wsButtons.Activate()
$wsButtons.Range("G15").Select()

Write-Host "Edit complete"
